$d = $word.ActiveDocument
$d.Content.Find.Execute("77-76=1", $true, $false, $false, $false, $false, $true, 1, $false, "61-50=11", 2) | Out-Null
$d.Content.Find.Execute("70-39=31", $true, $false, $false, $false, $false, $true, 1, $false, "0+86=86", 2) | Out-Null
$d.Content.Find.Execute("41+46=87", $true, $false, $false, $false, $false, $true, 1, $false, "60-28=32", 2) | Out-Null
$d.Content.Find.Execute("62-18=44", $true, $false, $false, $false, $false, $true, 1, $false, "58+11=69", 2) | Out-Null
$d.Content.Find.Execute("90-15=75", $true, $false, $false, $false, $false, $true, 1, $false, "43+6=49", 2) | Out-Null
$d.Content.Find.Execute("25-5=20", $true, $false, $false, $false, $false, $true, 1, $false, "64-33=31", 2) | Out-Null
$d.Content.Find.Execute("52+15=67", $true, $false, $false, $false, $false, $true, 1, $false, "99-8=91", 2) | Out-Null
$d.Content.Find.Execute("29+18=47", $true, $false, $false, $false, $false, $true, 1, $false, "91+2=93", 2) | Out-Null
$d.Content.Find.Execute("71+26=97", $true, $false, $false, $false, $false, $true, 1, $false, "1+24=25", 2) | Out-Null
$d.Content.Find.Execute("10+57=67", $true, $false, $false, $false, $false, $true, 1, $false, "71-20=51", 2) | Out-Null
$d.Content.Find.Execute("26+20=46", $true, $false, $false, $false, $false, $true, 1, $false, "50+37=87", 2) | Out-Null
$d.Content.Find.Execute("91-69=22", $true, $false, $false, $false, $false, $true, 1, $false, "76-43=33", 2) | Out-Null
$d.Content.Find.Execute("90-25=65", $true, $false, $false, $false, $false, $true, 1, $false, "37+45=82", 2) | Out-Null
$d.Content.Find.Execute("48+3=51", $true, $false, $false, $false, $false, $true, 1, $false, "85-28=57", 2) | Out-Null
$d.Content.Find.Execute("57-16=41", $true, $false, $false, $false, $false, $true, 1, $false, "20+8=28", 2) | Out-Null
$d.Content.Find.Execute("95-85=10", $true, $false, $false, $false, $false, $true, 1, $false, "71-62=9", 2) | Out-Null
$d.Content.Find.Execute("64-54=10", $true, $false, $false, $false, $false, $true, 1, $false, "48-19=29", 2) | Out-Null
$d.Content.Find.Execute("13+55=68", $true, $false, $false, $false, $false, $true, 1, $false, "49-7=42", 2) | Out-Null
$d.Content.Find.Execute("71-18=53", $true, $false, $false, $false, $false, $true, 1, $false, "14+82=96", 2) | Out-Null
$d.Content.Find.Execute("84-13=71", $true, $false, $false, $false, $false, $true, 1, $false, "85-25=60", 2) | Out-Null
$d.Content.Find.Execute("63+35=98", $true, $false, $false, $false, $false, $true, 1, $false, "31-17=14", 2) | Out-Null
$d.Content.Find.Execute("55+39=94", $true, $false, $false, $false, $false, $true, 1, $false, "38-21=17", 2) | Out-Null
$d.Content.Find.Execute("44+36=80", $true, $false, $false, $false, $false, $true, 1, $false, "26+63=89", 2) | Out-Null
$d.Content.Find.Execute("3+2=5", $true, $false, $false, $false, $false, $true, 1, $false, "68-42=26", 2) | Out-Null
$d.Content.Find.Execute("20+29=49", $true, $false, $false, $false, $false, $true, 1, $false, "14+73=87", 2) | Out-Null
$d.Content.Find.Execute("57-0=57", $true, $false, $false, $false, $false, $true, 1, $false, "72+12=84", 2) | Out-Null
$d.Content.Find.Execute("30-13=17", $true, $false, $false, $false, $false, $true, 1, $false, "93-50=43", 2) | Out-Null
$d.Content.Find.Execute("19+26=45", $true, $false, $false, $false, $false, $true, 1, $false, "27+51=78", 2) | Out-Null
$d.Content.Find.Execute("66-6=60", $true, $false, $false, $false, $false, $true, 1, $false, "71-8=63", 2) | Out-Null
$d.Content.Find.Execute("97-62=35", $true, $false, $false, $false, $false, $true, 1, $false, "68-54=14", 2) | Out-Null
$d.Content.Find.Execute("72-1=71", $true, $false, $false, $false, $false, $true, 1, $false, "54-39=15", 2) | Out-Null
$d.Content.Find.Execute("38-25=13", $true, $false, $false, $false, $false, $true, 1, $false, "44-28=16", 2) | Out-Null
$d.Content.Find.Execute("25+32=57", $true, $false, $false, $false, $false, $true, 1, $false, "87-87=0", 2) | Out-Null
$d.Content.Find.Execute("61+20=81", $true, $false, $false, $false, $false, $true, 1, $false, "34-24=10", 2) | Out-Null
$d.Content.Find.Execute("87-64=23", $true, $false, $false, $false, $false, $true, 1, $false, "0+63=63", 2) | Out-Null
$d.Content.Find.Execute("88-75=13", $true, $false, $false, $false, $false, $true, 1, $false, "0+85=85", 2) | Out-Null
$d.Content.Find.Execute("93-31=62", $true, $false, $false, $false, $false, $true, 1, $false, "52+45=97", 2) | Out-Null
$d.Content.Find.Execute("22+19=41", $true, $false, $false, $false, $false, $true, 1, $false, "60+29=89", 2) | Out-Null
$d.Content.Find.Execute("8+60=68", $true, $false, $false, $false, $false, $true, 1, $false, "36+42=78", 2) | Out-Null
$d.Content.Find.Execute("58+30=88", $true, $false, $false, $false, $false, $true, 1, $false, "58-45=13", 2) | Out-Null
$d.Content.Find.Execute("35+57=92", $true, $false, $false, $false, $false, $true, 1, $false, "54+20=74", 2) | Out-Null
$d.Content.Find.Execute("67-14=53", $true, $false, $false, $false, $false, $true, 1, $false, "8+19=27", 2) | Out-Null
$d.Content.Find.Execute("35-4=31", $true, $false, $false, $false, $false, $true, 1, $false, "18-10=8", 2) | Out-Null
$d.Content.Find.Execute("26-21=5", $true, $false, $false, $false, $false, $true, 1, $false, "31+28=59", 2) | Out-Null
$d.Content.Find.Execute("84-62=22", $true, $false, $false, $false, $false, $true, 1, $false, "95-44=51", 2) | Out-Null
$d.Content.Find.Execute("63-33=30", $true, $false, $false, $false, $false, $true, 1, $false, "99-97=2", 2) | Out-Null
$d.Content.Find.Execute("30+9=39", $true, $false, $false, $false, $false, $true, 1, $false, "46+45=91", 2) | Out-Null
$d.Content.Find.Execute("18+6=24", $true, $false, $false, $false, $false, $true, 1, $false, "36+44=80", 2) | Out-Null
$d.Content.Find.Execute("73-36=37", $true, $false, $false, $false, $false, $true, 1, $false, "24+72=96", 2) | Out-Null
$d.Content.Find.Execute("26+0=26", $true, $false, $false, $false, $false, $true, 1, $false, "86-64=22", 2) | Out-Null
$d.Content.Find.Execute("0+47=47", $true, $false, $false, $false, $false, $true, 1, $false, "58-29=29", 2) | Out-Null
$d.Content.Find.Execute("4+10=14", $true, $false, $false, $false, $false, $true, 1, $false, "77-52=25", 2) | Out-Null
$d.Content.Find.Execute("86-75=11", $true, $false, $false, $false, $false, $true, 1, $false, "59+14=73", 2) | Out-Null
$d.Content.Find.Execute("47-24=23", $true, $false, $false, $false, $false, $true, 1, $false, "36+33=69", 2) | Out-Null
$d.Content.Find.Execute("28+21=49", $true, $false, $false, $false, $false, $true, 1, $false, "15-9=6", 2) | Out-Null
$d.Content.Find.Execute("81-4=77", $true, $false, $false, $false, $false, $true, 1, $false, "91-0=91", 2) | Out-Null
$d.Content.Find.Execute("43-22=21", $true, $false, $false, $false, $false, $true, 1, $false, "15-6=9", 2) | Out-Null
$d.Content.Find.Execute("24+57=81", $true, $false, $false, $false, $false, $true, 1, $false, "96-37=59", 2) | Out-Null
$d.Content.Find.Execute("51-2=49", $true, $false, $false, $false, $false, $true, 1, $false, "62+1=63", 2) | Out-Null
$d.Content.Find.Execute("80-73=7", $true, $false, $false, $false, $false, $true, 1, $false, "36+10=46", 2) | Out-Null
$d.Content.Find.Execute("64-46=18", $true, $false, $false, $false, $false, $true, 1, $false, "79-17=62", 2) | Out-Null
$d.Content.Find.Execute("12+24=36", $true, $false, $false, $false, $false, $true, 1, $false, "43+44=87", 2) | Out-Null
$d.Content.Find.Execute("73-9=64", $true, $false, $false, $false, $false, $true, 1, $false, "55+9=64", 2) | Out-Null
$d.Content.Find.Execute("87-41=46", $true, $false, $false, $false, $false, $true, 1, $false, "59-24=35", 2) | Out-Null
$d.Content.Find.Execute("79-27=52", $true, $false, $false, $false, $false, $true, 1, $false, "12+16=28", 2) | Out-Null
$d.Content.Find.Execute("81-80=1", $true, $false, $false, $false, $false, $true, 1, $false, "82-27=55", 2) | Out-Null
$d.Content.Find.Execute("22+12=34", $true, $false, $false, $false, $false, $true, 1, $false, "88-7=81", 2) | Out-Null
$d.Content.Find.Execute("97-52=45", $true, $false, $false, $false, $false, $true, 1, $false, "19-2=17", 2) | Out-Null
$d.Content.Find.Execute("78-6=72", $true, $false, $false, $false, $false, $true, 1, $false, "18-9=9", 2) | Out-Null
$d.Content.Find.Execute("32+37=69", $true, $false, $false, $false, $false, $true, 1, $false, "6+16=22", 2) | Out-Null
$d.Content.Find.Execute("20+66=86", $true, $false, $false, $false, $false, $true, 1, $false, "55+0=55", 2) | Out-Null
$d.Content.Find.Execute("26+32=58", $true, $false, $false, $false, $false, $true, 1, $false, "52-40=12", 2) | Out-Null
$d.Content.Find.Execute("13+81=94", $true, $false, $false, $false, $false, $true, 1, $false, "82-24=58", 2) | Out-Null
$d.Content.Find.Execute("16+23=39", $true, $false, $false, $false, $false, $true, 1, $false, "55+22=77", 2) | Out-Null
$d.Content.Find.Execute("38+14=52", $true, $false, $false, $false, $false, $true, 1, $false, "19+23=42", 2) | Out-Null
$d.Content.Find.Execute("82-72=10", $true, $false, $false, $false, $false, $true, 1, $false, "10+76=86", 2) | Out-Null
$d.Content.Find.Execute("94-71=23", $true, $false, $false, $false, $false, $true, 1, $false, "79-33=46", 2) | Out-Null
$d.Content.Find.Execute("23+65=88", $true, $false, $false, $false, $false, $true, 1, $false, "86-1=85", 2) | Out-Null
$d.Content.Find.Execute("81-36=45", $true, $false, $false, $false, $false, $true, 1, $false, "74-60=14", 2) | Out-Null
$d.Content.Find.Execute("13+37=50", $true, $false, $false, $false, $false, $true, 1, $false, "5+86=91", 2) | Out-Null
$d.Content.Find.Execute("84+11=95", $true, $false, $false, $false, $false, $true, 1, $false, "33+60=93", 2) | Out-Null
$d.Content.Find.Execute("5+42=47", $true, $false, $false, $false, $false, $true, 1, $false, "89-29=60", 2) | Out-Null
$d.Content.Find.Execute("58+4=62", $true, $false, $false, $false, $false, $true, 1, $false, "19+69=88", 2) | Out-Null
$d.Content.Find.Execute("53-37=16", $true, $false, $false, $false, $false, $true, 1, $false, "86+0=86", 2) | Out-Null
$d.Content.Find.Execute("24+39=63", $true, $false, $false, $false, $false, $true, 1, $false, "89-10=79", 2) | Out-Null
$d.Content.Find.Execute("13+41=54", $true, $false, $false, $false, $false, $true, 1, $false, "67-65=2", 2) | Out-Null
$d.Content.Find.Execute("34+7=41", $true, $false, $false, $false, $false, $true, 1, $false, "49+0=49", 2) | Out-Null
$d.Content.Find.Execute("15+14=29", $true, $false, $false, $false, $false, $true, 1, $false, "34-0=34", 2) | Out-Null
$d.Content.Find.Execute("0+6=6", $true, $false, $false, $false, $false, $true, 1, $false, "67-63=4", 2) | Out-Null
$d.Content.Find.Execute("34+36=70", $true, $false, $false, $false, $false, $true, 1, $false, "65-42=23", 2) | Out-Null
$d.Content.Find.Execute("81-40=41", $true, $false, $false, $false, $false, $true, 1, $false, "94-60=34", 2) | Out-Null
$d.Content.Find.Execute("92-3=89", $true, $false, $false, $false, $false, $true, 1, $false, "51+0=51", 2) | Out-Null
$d.Content.Find.Execute("76-27=49", $true, $false, $false, $false, $false, $true, 1, $false, "58+25=83", 2) | Out-Null
$d.Content.Find.Execute("43-20=23", $true, $false, $false, $false, $false, $true, 1, $false, "35+41=76", 2) | Out-Null
$d.Content.Find.Execute("69-52=17", $true, $false, $false, $false, $false, $true, 1, $false, "73-69=4", 2) | Out-Null
$d.Content.Find.Execute("6+62=68", $true, $false, $false, $false, $false, $true, 1, $false, "4+61=65", 2) | Out-Null
$d.Content.Find.Execute("17+57=74", $true, $false, $false, $false, $false, $true, 1, $false, "21+20=41", 2) | Out-Null
$d.Content.Find.Execute("53-26=27", $true, $false, $false, $false, $false, $true, 1, $false, "59+30=89", 2) | Out-Null
$d.Content.Find.Execute("96-32=64", $true, $false, $false, $false, $false, $true, 1, $false, "75-68=7", 2) | Out-Null
$d.Content.Find.Execute("15+33=48", $true, $false, $false, $false, $false, $true, 1, $false, "84+14=98", 2) | Out-Null
